# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list ... with GitHub Actions".
#
# All cells in columns B:E on the sheet are plain text (inline/shared
# strings) in the source workbook, even though many Price values in column D
# look numeric (e.g. "343.63"). To avoid Excel silently re-typing those as
# real numbers (which would also reformat/round values such as "102.30" or
# "0.00000000362"), numeric-looking column D values are written with a
# leading apostrophe, Excel's standard force-text entry prefix; it is not
# stored as part of the cell value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.703.08"
$ws.Range("E2").Value = "  -1.25%  "
$ws.Range("D3").Value = "'2.098.37"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  +0.45%  "
$ws.Range("D5").Value = "'343.63"
$ws.Range("E5").Value = "  -1.77%  "
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("D7").Value = "'0.5171"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "'0.4378"
$ws.Range("E8").Value = "  -2.22%  "
$ws.Range("D9").Value = "'53.34"
$ws.Range("E9").Value = "  +1.29%  "
$ws.Range("D10").Value = "'0.09199"
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("D11").Value = "'1.167"
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("D12").Value = "'24.63"
$ws.Range("E12").Value = "  -4.54%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'2.080.00"
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'6.768"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").Value = "'8.146"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").Value = "'102.30"
$ws.Range("E16").Value = "  +3.16%  "
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").Value = "'21.05"
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("D20").Value = "'0.06666"
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("D22").Value = "'6.201"
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("D23").Value = "'29.750.19"
$ws.Range("E23").Value = "  -1.42%  "
$ws.Range("D24").Value = "'12.61"
$ws.Range("E24").Value = "  -2.33%  "
$ws.Range("D25").Value = "'2.304"
$ws.Range("E25").Value = "  -2.12%  "
$ws.Range("D26").Value = "'2.306.22"
$ws.Range("E26").Value = "  -2.31%  "
$ws.Range("D27").Value = "'21.89"
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("D29").Value = "'2.495"
$ws.Range("E29").Value = "  -2.55%  "
$ws.Range("D30").Value = "'133.42"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "'1.129"
$ws.Range("E31").Value = "  -4.33%  "
$ws.Range("D32").Value = "'1.682"
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("E33").Value = "  -1.75%  "
$ws.Range("D34").Value = "'6.193"
$ws.Range("E34").Value = "  -1.25%  "
$ws.Range("D35").Value = "'3.959"
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("D36").Value = "'6.329"
$ws.Range("E36").Value = "  +6.78%  "
$ws.Range("D37").Value = "'10.42"
$ws.Range("E37").Value = "  +2.11%  "
$ws.Range("D38").Value = "'0.02579"
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("D39").Value = "'0.06707"
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("D40").Value = "'0.6986"
$ws.Range("E40").Value = "  +2.13%  "
$ws.Range("D41").Value = "'1.331"
$ws.Range("E41").Value = "  +5.99%  "
$ws.Range("D42").Value = "'12.42"
$ws.Range("E42").Value = "  -3.34%  "
$ws.Range("D43").Value = "'0.2212"
$ws.Range("E43").Value = "  -4.44%  "
$ws.Range("D44").Value = "'0.6797"
$ws.Range("E44").Value = "  +5.82%  "
$ws.Range("D45").Value = "'14.27"
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("D46").Value = "'2.319"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Value = "'0.00000000362"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").Value = "'3.611"
$ws.Range("E48").Value = "  -1.41%  "
$ws.Range("D50").Value = "'1.216"
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("D51").Value = "'81.14"
$ws.Range("E51").Value = "  -2.99%  "
